# Trace_Report_NOVAMEAL_initial.xlsx — refresh the single trace-event row
# with a newer search result (new location/event/time/weights) and strip
# the leftover AutoFilter / _FilterDatabase defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds the free-text search-description banner; only the timestamp
# embedded in it changed for this re-run of the trace search.
$ws.Range("A1").Value = "Description unknown, completed 06/15/2023 05:57:14 EDT, by WPJTOWN1.The search returned: 1 events."

# Row 2 is the header row and is unchanged.

# Row 3 is the single data row — update it to the new event details.
$ws.Range("C3").Value = "JOHNSTOWN"   # Location City
$ws.Range("D3").Value = "CO"          # State
$ws.Range("E3").Value = 6             # Month
$ws.Range("F3").Value = 1             # Day
$ws.Range("G3").Value = 1811          # Time
$ws.Range("H3").Value = "Placed Actual" # Event
$ws.Range("I3").ClearContents()       # Train ID no longer reported - clear it
$ws.Range("J3").Value = "LOVELAND"    # Destination City
$ws.Range("K3").Value = "CO"          # State (destination)

# Gross/Tare/Net weight (L3:N3) and Car_no (O3) are unchanged.

# The AutoFilter (and its accompanying hidden _FilterDatabase defined name)
# is no longer present on the refreshed sheet.
$ws.AutoFilterMode = $false
foreach ($n in $wb.Names) {
    $n.Delete()
}
